$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-5 from 45243 to 45244
$ws.Range("C2:C5").Value = 45244
